$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44335
$ws.Range("N2").Value = 19500
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19750
$ws.Range("S2").Value = 1097
$ws.Range("D3").Value = 44335
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 17500
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17750
$ws.Range("S3").Value = 986
$ws.Range("D4").Value = 44335
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 12500
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12750
$ws.Range("S4").Value = 708
$ws.Range("D5").Value = 44294
$ws.Range("M5").Value = 360
$ws.Range("D6").Value = 44294
$ws.Range("M6").Value = 240
$ws.Range("D7").Value = 44294
$ws.Range("M7").Value = 240
$ws.Range("D8").Value = 44384
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 14500
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14750
$ws.Range("S8").Value = 819
$ws.Range("D9").Value = 44384
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = 11500
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 11750
$ws.Range("S9").Value = 653
$ws.Range("D10").Value = 44384
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 8500
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 8750
$ws.Range("S10").Value = 486
$ws.Range("D11").Value = 44384
$ws.Range("L11").Value = "Tercera"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 5500
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 5750
$ws.Range("S11").Value = 319
$ws.Range("D12").Value = 44224
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 16500
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16750
$ws.Range("S12").Value = 931
$ws.Range("D13").Value = 44224
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 13500
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 13750
$ws.Range("S13").Value = 764
$ws.Range("D14").Value = 44272
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 12500
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 12750
$ws.Range("S14").Value = 708
$ws.Range("D15").Value = 44272
$ws.Range("N15").Value = 10500
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 10750
$ws.Range("S15").Value = 597
$ws.Range("D16").Value = 44272
$ws.Range("N16").Value = 8500
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 8750
$ws.Range("S16").Value = 486
$ws.Range("D17").Value = 44293
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 12500
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 12750
$ws.Range("S17").Value = 708
$ws.Range("D18").Value = 44293
$ws.Range("M18").Value = 508
$ws.Range("N18").Value = 10500
$ws.Range("O18").Value = 11000
$ws.Range("P18").Value = 10746
$ws.Range("S18").Value = 597
$ws.Range("D19").Value = 44293
$ws.Range("M19").Value = 400
$ws.Range("N19").Value = 8500
$ws.Range("O19").Value = 9000
$ws.Range("P19").Value = 8750
$ws.Range("S19").Value = 486
$ws.Range("D20").Value = 44385
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14500
$ws.Range("P20").Value = 14250
$ws.Range("S20").Value = 792
$ws.Range("D21").Value = 44385
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 11500
$ws.Range("P21").Value = 11250
$ws.Range("S21").Value = 625
$ws.Range("D22").Value = 44385
$ws.Range("M22").Value = 240
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8500
$ws.Range("P22").Value = 8250
$ws.Range("S22").Value = 458
$ws.Range("D23").Value = 44385
$ws.Range("M23").Value = 120
$ws.Range("N23").Value = 5000
$ws.Range("O23").Value = 5500
$ws.Range("P23").Value = 5250
$ws.Range("S23").Value = 292
$ws.Range("D24").Value = 44280
$ws.Range("M24").Value = 240
$ws.Range("D25").Value = 44280
$ws.Range("M25").Value = 240
$ws.Range("D26").Value = 44280
$ws.Range("M26").Value = 300
$ws.Range("D27").Value = 44308
$ws.Range("M27").Value = 300
$ws.Range("N27").Value = 15500
$ws.Range("O27").Value = 16000
$ws.Range("P27").Value = 15750
$ws.Range("S27").Value = 875
$ws.Range("D28").Value = 44308
$ws.Range("M28").Value = 240
$ws.Range("N28").Value = 13500
$ws.Range("O28").Value = 14000
$ws.Range("P28").Value = 13750
$ws.Range("S28").Value = 764
$ws.Range("D29").Value = 44308
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 10500
$ws.Range("O29").Value = 11000
$ws.Range("P29").Value = 10750
$ws.Range("S29").Value = 597
$ws.Range("D30").Value = 44279
$ws.Range("M30").Value = 200
$ws.Range("D31").Value = 44279
$ws.Range("D32").Value = 44279
$ws.Range("M32").Value = 240
$ws.Range("D33").Value = 44273
$ws.Range("M33").Value = 160
$ws.Range("D34").Value = 44273
$ws.Range("D35").Value = 44273
$ws.Range("M35").Value = 200
$ws.Range("D36").Value = 44286
$ws.Range("M36").Value = 700
$ws.Range("D37").Value = 44286
$ws.Range("M37").Value = 500
$ws.Range("D38").Value = 44286
$ws.Range("M38").Value = 300
